# Add a new "assignments" worksheet right after "attendance", mirroring the
# ID/Name columns from attendance and recording each student's assignment
# marks (out of 10), with Percentage and Grade formulas just like the other
# gradebook sheets (quiz1, quiz2, ...).

$wb = $excel.ActiveWorkbook

$attendance = $wb.Worksheets.Item("attendance")

$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $attendance)
$ws.Name = "assignments"

# Headers
$ws.Range("A1").Value2 = "ID"
$ws.Range("B1").Value2 = "Name"
$ws.Range("C1").Value2 = "Obtained"
$ws.Range("D1").Value2 = "Percentage"
$ws.Range("E1").Value2 = "Grade"

# Student rows: ID, Name, Marks obtained (out of 10)
$students = @(
  @(183011218, "Aidid Rashed Efat", 9),
  @(173014033, "Mohammad Moniruzzaman Mollah", 0),
  @(181014001, "*Towfiq Imroze", 0),
  @(181014051, "*Emtiaz Hossain Tamim", 0),
  @(181014126, "*Md. Mehedi Hasan Neloy", 8),
  @(182014056, "*Dipto Kumar Pramanik", 8),
  @(183014070, "*Mahmudul Hasan", 0),
  @(191014032, "*Ashraf Uddin Tushar", 0),
  @(193014009, "*Anika Tabassum", 10),
  @(193014067, "*Asib Sikder Pranto", 8),
  @(193014071, "*Jannatul Nayem", 10),
  @(193014072, "*G.M. Shahriar Rahman", 0),
  @(193014073, "*Asjad Hossain Khan", 10),
  @(183016002, "*Farhan Bin Murtaza", 0),
  @(192016001, "*Joy Saha", 9),
  @(193016008, "*Md. Shojal Hossain", 8)
)

for ($i = 0; $i -lt $students.Count; $i++) {
  $r = 2 + $i
  $ws.Cells.Item($r, 1).Value2 = $students[$i][0]
  $ws.Cells.Item($r, 2).Value2 = $students[$i][1]
  $ws.Cells.Item($r, 3).Value2 = $students[$i][2]
  $ws.Cells.Item($r, 4).Formula = "=(C$r/C`$18)*100"
  $ws.Cells.Item($r, 5).Formula = "=IF(D$r>94,""A+"",IF(D$r>84,""A"",IF(D$r>79,""A-"",IF(D$r>74,""B+"",IF(D$r>69,""B"",IF(D$r>64,""B-"",IF(D$r>59,""C+"",IF(D$r>54,""C"",IF(D$r>49,""D"",""F"")))))))))"
}

# Row 18 holds the total possible marks for the assignments (10)
$ws.Range("C18").Value2 = 10

# Column widths to match the rest of the workbook
$ws.Columns.Item(1).ColumnWidth = 14.6640625
$ws.Columns.Item(2).ColumnWidth = 30.6640625
$ws.Columns.Item(3).ColumnWidth = 14.6640625
$ws.Columns.Item(4).ColumnWidth = 14.6640625
$ws.Columns.Item(5).ColumnWidth = 14.6640625

# Conditional formatting: highlight "F" grades in column E (rows 2-17)
$rngF = $ws.Range("E2:E17")
$cf = $rngF.FormatConditions.Add(1, 3, '"F"')
$cf.Interior.Color = 13551615
$cf.Font.Color = 192

$ws.Range("A1").Select()
